$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.013.78"
$ws.Range("E2").Value = "  -1.77%  "

$ws.Range("D3").Value = "1.829.93"
$ws.Range("E3").Value = "  -1.02%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.31%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.37"
$ws.Range("E5").Value = "  -2.81%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.30%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4654"
$ws.Range("E7").Value = "  -0.27%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3861"
$ws.Range("E8").Value = "  -1.57%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07860"
$ws.Range("E9").Value = "  -0.62%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9587"
$ws.Range("E10").Value = "  -2.64%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.86"
$ws.Range("E11").Value = "  -1.51%  "

$ws.Range("D12").Value = "2.000.90"
$ws.Range("E12").Value = "  +6.25%  "

$ws.Range("E13").Value = "  -3.05%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.892"
$ws.Range("E14").Value = "  -1.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06829"
$ws.Range("E15").Value = "  -0.35%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.23"
$ws.Range("E16").Value = "  -0.56%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("E17").Value = "  -0.35%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009915"
$ws.Range("E18").Value = "  -1.48%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.56"
$ws.Range("E19").Value = "  -2.96%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  -0.26%  "

$ws.Range("D21").Value = "28.005.00"
$ws.Range("E21").Value = "  -1.89%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.314"
$ws.Range("E22").Value = "  -1.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.98"
$ws.Range("E23").Value = "  -2.38%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.088"
$ws.Range("E24").Value = "  -1.83%  "

$ws.Range("D25").Value = "2.036.07"
$ws.Range("E25").Value = "  -6.18%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.60"
$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.07"
$ws.Range("E27").Value = "  -1.64%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.720"
$ws.Range("E28").Value = "  -6.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.959"
$ws.Range("E29").Value = "  -2.78%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.51"
$ws.Range("E30").Value = "  +0.14%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.9350"
$ws.Range("E31").Value = "  -4.37%  "

$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09237"
$ws.Range("E32").Value = "  -1.92%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.267"
$ws.Range("E33").Value = "  -1.80%  "

$ws.Range("E34").Value = "  -2.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.292"
$ws.Range("E35").Value = "  -5.95%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05857"
$ws.Range("E36").Value = "  -4.37%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02141"
$ws.Range("E37").Value = "  -2.45%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.139"
$ws.Range("E38").Value = "  -1.77%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.789"
$ws.Range("E39").Value = "  +2.62%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5574"
$ws.Range("E40").Value = "  -2.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.840"
$ws.Range("E41").Value = "  -2.47%  "

$ws.Range("E42").Value = "  -1.97%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.59"
$ws.Range("E43").Value = "  -1.80%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.07020"
$ws.Range("E44").Value = "  -1.92%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5241"
$ws.Range("E45").Value = "  -2.40%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.121"
$ws.Range("E46").Value = "  -11.43%  "

$ws.Range("E47").Value = "  -4.17%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "112.76"
$ws.Range("E48").Value = "  -0.77%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.102"
$ws.Range("E49").Value = "  -9.98%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.9996"
$ws.Range("E50").Value = "  -0.28%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.320"
